$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 548, shifting the existing
# rows 548:566 down to 549:567 (all of their values move with them).
$ws.Rows.Item(548).EntireRow.Insert()

# Populate the newly inserted row 548 with a new price record. Most
# columns mirror the row that used to occupy this position (now row
# 549) - only Fecha (D), Volumen (M) and Origen (R) differ.
$ws.Cells.Item(548, 1).Value = 10
$ws.Cells.Item(548, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(548, 3).Value = "La Araucanía"
$ws.Cells.Item(548, 4).Value = 45075
$ws.Cells.Item(548, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(548, 5).Value = 9
$ws.Cells.Item(548, 6).Value = "Fruta"
$ws.Cells.Item(548, 7).Value = 100108
$ws.Cells.Item(548, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(548, 9).Value = 100108002
$ws.Cells.Item(548, 10).Value = "Mango"
$ws.Cells.Item(548, 11).Value = "Sin especificar"
$ws.Cells.Item(548, 12).Value = "Primera"
$ws.Cells.Item(548, 13).Value = 1800
$ws.Cells.Item(548, 14).Value = 9000
$ws.Cells.Item(548, 15).Value = 9000
$ws.Cells.Item(548, 16).Value = 9000
$ws.Cells.Item(548, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(548, 18).Value = "Perú"
$ws.Cells.Item(548, 19).Value = 2250
$ws.Cells.Item(548, 20).Value = 4
